$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 330
$ws1.Range("F5").Value = 182
$ws1.Range("F6").Value = 687
$ws1.Range("F8").Value = 492
$ws1.Range("F9").Value = 89
$ws1.Range("F10").Value = 531
$ws1.Range("F11").Value = 432
$ws1.Range("F12").Value = 69
$ws1.Range("F13").Value = 33
$ws1.Range("F14").Value = 120
$ws1.Range("F15").Value = 205

# Sheet "本地生活" (sheet3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6233
$ws3.Range("F4").Value = 758
$ws3.Range("F5").Value = 1838

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6233
$ws4.Range("F4").Value = 758
$ws4.Range("F5").Value = 1838
$ws4.Range("F6").Value = 330
$ws4.Range("F12").Value = 182
$ws4.Range("F15").Value = 687
$ws4.Range("F19").Value = 492
$ws4.Range("F21").Value = 89
$ws4.Range("F22").Value = 531
$ws4.Range("F24").Value = 432
$ws4.Range("F25").Value = 69
$ws4.Range("F28").Value = 33
$ws4.Range("F29").Value = 120
$ws4.Range("F35").Value = 205
